$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "PARTNER & STRATEGIC CONSULTANT - Siege Analytics, Washington, DC | January 2014 – Present" "PARTNER - Siege Analytics, Washington, DC | January 2014 – Present"
Replace-Text "PRINCIPAL MARKETING CONSULTANT - Clarity and Rigour, Washington, DC | 2012 – 2014" "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014"
Replace-Text "DIRECTOR OF MARKETING - Helm, Washington, DC | 2010 – 2012" "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012"
Replace-Text "SENIOR MARKETING ANALYST - GSD&M, Austin, TX | 2008 – 2010" "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010"
Replace-Text "MARKETING COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008" "RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008"
Replace-Text "MARKETING SPECIALIST - Salsa Labs, Inc., Washington, DC | 2004 – 2006" "SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006"
Replace-Text "COMMUNICATIONS COORDINATOR - The Praxis Project, Oakland, CA | 2002 – 2004" "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004"
Replace-Text "RESEARCH COORDINATOR - Lake Research Partners, Washington, DC | 2001 – 2002" "PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002"
Replace-Text "FIELD COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001" "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001"
